$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Team name swaps (plain text, direct assignment) ---
$ws.Range("B5").Value = "Bologna"
$ws.Range("B6").Value = "Roma"
$ws.Range("B7").Value = "Atalanta"
$ws.Range("B8").Value = "Napoli"
$ws.Range("B10").Value = "Lazio"
$ws.Range("B15").Value = "Empoli"
$ws.Range("B16").Value = "Udinese"
$ws.Range("B17").Value = "Frosinone"
$ws.Range("B19").Value = "Cagliari"
$ws.Range("B20").Value = "Sassuolo"

# --- Numeric-looking values (odds %, decimals) must be forced to text
# so they stay as shared strings like the rest of the table, matching
# the original "1.4" / "85%" / "2.96" style text cells instead of being
# auto-converted to real numbers by Excel. We stage each value in a
# scratch cell formatted as Text, copy it, and paste-special (values only)
# into the real target cell, then clear the scratch cell.
$scratch = $ws.Range("Z1")
$scratch.NumberFormat = "@"

$scratch.Value = "1.4"
$scratch.Copy()
$ws.Range("C2").PasteSpecial(-4163)
$scratch.Value = "85%"
$scratch.Copy()
$ws.Range("E2").PasteSpecial(-4163)
$scratch.Value = "59%"
$scratch.Copy()
$ws.Range("F2").PasteSpecial(-4163)
$scratch.Value = "3.04"
$scratch.Copy()
$ws.Range("G2").PasteSpecial(-4163)
$scratch.Value = "2.4"
$scratch.Copy()
$ws.Range("C3").PasteSpecial(-4163)
$scratch.Value = "5.0"
$scratch.Copy()
$ws.Range("D3").PasteSpecial(-4163)
$scratch.Value = "66%"
$scratch.Copy()
$ws.Range("E3").PasteSpecial(-4163)
$scratch.Value = "44%"
$scratch.Copy()
$ws.Range("F3").PasteSpecial(-4163)
$scratch.Value = "2.33"
$scratch.Copy()
$ws.Range("G3").PasteSpecial(-4163)
$scratch.Value = "2.4"
$scratch.Copy()
$ws.Range("C4").PasteSpecial(-4163)
$scratch.Value = "70%"
$scratch.Copy()
$ws.Range("E4").PasteSpecial(-4163)
$scratch.Value = "59%"
$scratch.Copy()
$ws.Range("F4").PasteSpecial(-4163)
$scratch.Value = "3.07"
$scratch.Copy()
$ws.Range("G4").PasteSpecial(-4163)
$scratch.Value = "2.4"
$scratch.Copy()
$ws.Range("C5").PasteSpecial(-4163)
$scratch.Value = "4.0"
$scratch.Copy()
$ws.Range("D5").PasteSpecial(-4163)
$scratch.Value = "82%"
$scratch.Copy()
$ws.Range("E5").PasteSpecial(-4163)
$scratch.Value = "49%"
$scratch.Copy()
$ws.Range("F5").PasteSpecial(-4163)
$scratch.Value = "2.41"
$scratch.Copy()
$ws.Range("G5").PasteSpecial(-4163)
$scratch.Value = "2.5"
$scratch.Copy()
$ws.Range("C6").PasteSpecial(-4163)
$scratch.Value = "4.3"
$scratch.Copy()
$ws.Range("D6").PasteSpecial(-4163)
$scratch.Value = "85%"
$scratch.Copy()
$ws.Range("E6").PasteSpecial(-4163)
$scratch.Value = "63%"
$scratch.Copy()
$ws.Range("F6").PasteSpecial(-4163)
$scratch.Value = "3.15"
$scratch.Copy()
$ws.Range("G6").PasteSpecial(-4163)
$scratch.Value = "2.3"
$scratch.Copy()
$ws.Range("C7").PasteSpecial(-4163)
$scratch.Value = "5.7"
$scratch.Copy()
$ws.Range("D7").PasteSpecial(-4163)
$scratch.Value = "85%"
$scratch.Copy()
$ws.Range("E7").PasteSpecial(-4163)
$scratch.Value = "59%"
$scratch.Copy()
$ws.Range("F7").PasteSpecial(-4163)
$scratch.Value = "2.93"
$scratch.Copy()
$ws.Range("G7").PasteSpecial(-4163)
$scratch.Value = "2.2"
$scratch.Copy()
$ws.Range("C8").PasteSpecial(-4163)
$scratch.Value = "6.0"
$scratch.Copy()
$ws.Range("D8").PasteSpecial(-4163)
$scratch.Value = "78%"
$scratch.Copy()
$ws.Range("E8").PasteSpecial(-4163)
$scratch.Value = "60%"
$scratch.Copy()
$ws.Range("F8").PasteSpecial(-4163)
$scratch.Value = "2.70"
$scratch.Copy()
$ws.Range("G8").PasteSpecial(-4163)
$scratch.Value = "5.0"
$scratch.Copy()
$ws.Range("D9").PasteSpecial(-4163)
$scratch.Value = "67%"
$scratch.Copy()
$ws.Range("E9").PasteSpecial(-4163)
$scratch.Value = "2.56"
$scratch.Copy()
$ws.Range("G9").PasteSpecial(-4163)
$scratch.Value = "2.6"
$scratch.Copy()
$ws.Range("C10").PasteSpecial(-4163)
$scratch.Value = "4.9"
$scratch.Copy()
$ws.Range("D10").PasteSpecial(-4163)
$scratch.Value = "70%"
$scratch.Copy()
$ws.Range("E10").PasteSpecial(-4163)
$scratch.Value = "40%"
$scratch.Copy()
$ws.Range("F10").PasteSpecial(-4163)
$scratch.Value = "2.26"
$scratch.Copy()
$ws.Range("G10").PasteSpecial(-4163)
$scratch.Value = "2.1"
$scratch.Copy()
$ws.Range("C11").PasteSpecial(-4163)
$scratch.Value = "63%"
$scratch.Copy()
$ws.Range("E11").PasteSpecial(-4163)
$scratch.Value = "30%"
$scratch.Copy()
$ws.Range("F11").PasteSpecial(-4163)
$scratch.Value = "1.85"
$scratch.Copy()
$ws.Range("G11").PasteSpecial(-4163)
$scratch.Value = "2.3"
$scratch.Copy()
$ws.Range("C12").PasteSpecial(-4163)
$scratch.Value = "5.3"
$scratch.Copy()
$ws.Range("D12").PasteSpecial(-4163)
$scratch.Value = "67%"
$scratch.Copy()
$ws.Range("E12").PasteSpecial(-4163)
$scratch.Value = "37%"
$scratch.Copy()
$ws.Range("F12").PasteSpecial(-4163)
$scratch.Value = "2.30"
$scratch.Copy()
$ws.Range("G12").PasteSpecial(-4163)
$scratch.Value = "4.0"
$scratch.Copy()
$ws.Range("D13").PasteSpecial(-4163)
$scratch.Value = "67%"
$scratch.Copy()
$ws.Range("E13").PasteSpecial(-4163)
$scratch.Value = "41%"
$scratch.Copy()
$ws.Range("F13").PasteSpecial(-4163)
$scratch.Value = "2.30"
$scratch.Copy()
$ws.Range("G13").PasteSpecial(-4163)
$scratch.Value = "2.6"
$scratch.Copy()
$ws.Range("C14").PasteSpecial(-4163)
$scratch.Value = "4.3"
$scratch.Copy()
$ws.Range("D14").PasteSpecial(-4163)
$scratch.Value = "82%"
$scratch.Copy()
$ws.Range("E14").PasteSpecial(-4163)
$scratch.Value = "45%"
$scratch.Copy()
$ws.Range("F14").PasteSpecial(-4163)
$scratch.Value = "2.56"
$scratch.Copy()
$ws.Range("G14").PasteSpecial(-4163)
$scratch.Value = "2.1"
$scratch.Copy()
$ws.Range("C15").PasteSpecial(-4163)
$scratch.Value = "4.9"
$scratch.Copy()
$ws.Range("D15").PasteSpecial(-4163)
$scratch.Value = "67%"
$scratch.Copy()
$ws.Range("E15").PasteSpecial(-4163)
$scratch.Value = "37%"
$scratch.Copy()
$ws.Range("F15").PasteSpecial(-4163)
$scratch.Value = "2.33"
$scratch.Copy()
$ws.Range("G15").PasteSpecial(-4163)
$scratch.Value = "2.4"
$scratch.Copy()
$ws.Range("C16").PasteSpecial(-4163)
$scratch.Value = "4.4"
$scratch.Copy()
$ws.Range("D16").PasteSpecial(-4163)
$scratch.Value = "77%"
$scratch.Copy()
$ws.Range("E16").PasteSpecial(-4163)
$scratch.Value = "40%"
$scratch.Copy()
$ws.Range("F16").PasteSpecial(-4163)
$scratch.Value = "2.48"
$scratch.Copy()
$ws.Range("G16").PasteSpecial(-4163)
$scratch.Value = "2.0"
$scratch.Copy()
$ws.Range("C17").PasteSpecial(-4163)
$scratch.Value = "5.4"
$scratch.Copy()
$ws.Range("D17").PasteSpecial(-4163)
$scratch.Value = "93%"
$scratch.Copy()
$ws.Range("E17").PasteSpecial(-4163)
$scratch.Value = "70%"
$scratch.Copy()
$ws.Range("F17").PasteSpecial(-4163)
$scratch.Value = "3.37"
$scratch.Copy()
$ws.Range("G17").PasteSpecial(-4163)
$scratch.Value = "2.3"
$scratch.Copy()
$ws.Range("C18").PasteSpecial(-4163)
$scratch.Value = "3.3"
$scratch.Copy()
$ws.Range("D18").PasteSpecial(-4163)
$scratch.Value = "60%"
$scratch.Copy()
$ws.Range("E18").PasteSpecial(-4163)
$scratch.Value = "45%"
$scratch.Copy()
$ws.Range("F18").PasteSpecial(-4163)
$scratch.Value = "2.22"
$scratch.Copy()
$ws.Range("G18").PasteSpecial(-4163)
$scratch.Value = "2.1"
$scratch.Copy()
$ws.Range("C19").PasteSpecial(-4163)
$scratch.Value = "4.7"
$scratch.Copy()
$ws.Range("D19").PasteSpecial(-4163)
$scratch.Value = "82%"
$scratch.Copy()
$ws.Range("E19").PasteSpecial(-4163)
$scratch.Value = "56%"
$scratch.Copy()
$ws.Range("F19").PasteSpecial(-4163)
$scratch.Value = "2.67"
$scratch.Copy()
$ws.Range("G19").PasteSpecial(-4163)
$scratch.Value = "1.8"
$scratch.Copy()
$ws.Range("C20").PasteSpecial(-4163)
$scratch.Value = "5.4"
$scratch.Copy()
$ws.Range("D20").PasteSpecial(-4163)
$scratch.Value = "82%"
$scratch.Copy()
$ws.Range("E20").PasteSpecial(-4163)
$scratch.Value = "59%"
$scratch.Copy()
$ws.Range("F20").PasteSpecial(-4163)
$scratch.Value = "3.22"
$scratch.Copy()
$ws.Range("G20").PasteSpecial(-4163)
$scratch.Value = "2.3"
$scratch.Copy()
$ws.Range("C21").PasteSpecial(-4163)
$scratch.Value = "4.2"
$scratch.Copy()
$ws.Range("D21").PasteSpecial(-4163)
$scratch.Value = "85%"
$scratch.Copy()
$ws.Range("E21").PasteSpecial(-4163)
$scratch.Value = "63%"
$scratch.Copy()
$ws.Range("F21").PasteSpecial(-4163)
$scratch.Value = "2.78"
$scratch.Copy()
$ws.Range("G21").PasteSpecial(-4163)

$scratch.Clear()
